$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "باسل العنزي"
$ws.Range("B69").Value = 500022150
$ws.Range("C69").Value = "الهفوف"

$ws.Range("A70").Value = "البراء السحيباني"
$ws.Range("B70").Value = "+966 53 412 8937"
$ws.Range("C70").Value = "العقير"

$ws.Range("C70").Select()

